$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $rowIndex, $newText) {
    $cell = $table.Cell($rowIndex, 1)
    $r = $cell.Range
    # Exclude the trailing cell-mark character so we don't overwrite it
    $r.End = $r.End - 1
    $r.Text = $newText
}

# Simple single-value cell replacements
Set-CellText $t 1 "0M"
Set-CellText $t 2 "0M"
Set-CellText $t 3 "0M"
Set-CellText $t 4 "47"
Set-CellText $t 6 "0.03542"
Set-CellText $t 7 "0.00510"
Set-CellText $t 8 "0.00000"
Set-CellText $t 9 "0.03542"
Set-CellText $t 10 "0.03542"
Set-CellText $t 11 "0.03542"
Set-CellText $t 12 "0.03738"

# Collapse the multi-run rows (44, 45, 46) into single values
Set-CellText $t 44 "99.94"
Set-CellText $t 45 "0.04"
Set-CellText $t 46 "61"
